# GuildConfig.xlsx edit
# Commit: "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to be a generic "Property1" table is renamed to
# "DataNode" to match the new DataNode/DataTable/Entity naming scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet.
$ws.Name = "DataNode"

# The author's selection/cursor ended up on D36 (frozen-pane bottom-left
# section) when the workbook was last saved.
[void]$ws.Range("D36").Select()
